# Refresh the cryptocurrency price/volume snapshot (Price column D, 1h Volume
# change column E) with the latest scraped values, as captured by the
# "Updated cryptos list" GitHub Actions job.
#
# Some Price values are plain decimals (e.g. "252.62", "5.00") that Excel's
# automatic type inference would otherwise coerce into numbers - silently
# dropping significant trailing zeros (5.00 -> 5) or reformatting the text.
# The source data treats this whole column as text (prices like
# "35.101.56" / "1.900.22" use '.' as a thousands separator, not a decimal
# point), so any cell whose new value looks like a plain number is forced
# to Text format before the value is written, keeping it an exact string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '35.101.56' },
    @{ Cell = 'E2'; Value = '  -0.55%  ' },
    @{ Cell = 'D3'; Value = '1.900.22' },
    @{ Cell = 'E3'; Value = '  +0.07%  ' },
    @{ Cell = 'E4'; Value = '  -0.47%  ' },
    @{ Cell = 'D5'; Value = '252.62' },
    @{ Cell = 'E5'; Value = '  +2.58%  ' },
    @{ Cell = 'E6'; Value = '  +1.15%  ' },
    @{ Cell = 'E7'; Value = '  -0.39%  ' },
    @{ Cell = 'E8'; Value = '  +3.10%  ' },
    @{ Cell = 'D9'; Value = '0.355' },
    @{ Cell = 'E9'; Value = '  +1.77%  ' },
    @{ Cell = 'D10'; Value = '52.29' },
    @{ Cell = 'E10'; Value = '  +0.66%  ' },
    @{ Cell = 'D11'; Value = '0.0758' },
    @{ Cell = 'E11'; Value = '  +5.00%  ' },
    @{ Cell = 'D12'; Value = '0.0977' },
    @{ Cell = 'E12'; Value = '  -1.01%  ' },
    @{ Cell = 'D13'; Value = '13.06' },
    @{ Cell = 'E13'; Value = '  +4.23%  ' },
    @{ Cell = 'D14'; Value = '2.176.37' },
    @{ Cell = 'E14'; Value = '  +0.11%  ' },
    @{ Cell = 'E15'; Value = '  +3.28%  ' },
    @{ Cell = 'D16'; Value = '5.00' },
    @{ Cell = 'E16'; Value = '  +3.33%  ' },
    @{ Cell = 'D17'; Value = '1.891.09' },
    @{ Cell = 'E17'; Value = '  -0.45%  ' },
    @{ Cell = 'D18'; Value = '35.139.57' },
    @{ Cell = 'E18'; Value = '  -0.40%  ' },
    @{ Cell = 'D19'; Value = '74.12' },
    @{ Cell = 'E19'; Value = '  +2.22%  ' },
    @{ Cell = 'D20'; Value = '0.0₃0840' },
    @{ Cell = 'E20'; Value = '  +2.66%  ' },
    @{ Cell = 'D21'; Value = '251.62' },
    @{ Cell = 'E21'; Value = '  +4.41%  ' },
    @{ Cell = 'E22'; Value = '  +1.81%  ' },
    @{ Cell = 'D23'; Value = '5.05' },
    @{ Cell = 'E23'; Value = '  +0.84%  ' },
    @{ Cell = 'E24'; Value = '  -0.40%  ' },
    @{ Cell = 'E25'; Value = '  +5.84%  ' },
    @{ Cell = 'D26'; Value = '2.28' },
    @{ Cell = 'E26'; Value = '  -1.47%  ' },
    @{ Cell = 'D27'; Value = '168.57' },
    @{ Cell = 'E27'; Value = '  +0.52%  ' },
    @{ Cell = 'D28'; Value = '8.60' },
    @{ Cell = 'E28'; Value = '  +0.35%  ' },
    @{ Cell = 'D29'; Value = '18.52' },
    @{ Cell = 'E29'; Value = '  -2.30%  ' },
    @{ Cell = 'E30'; Value = '  -0.72%  ' },
    @{ Cell = 'D31'; Value = '4.128.69' },
    @{ Cell = 'E31'; Value = '  -0.33%  ' },
    @{ Cell = 'D32'; Value = '2.07' },
    @{ Cell = 'E32'; Value = '  +10.77%  ' },
    @{ Cell = 'D33'; Value = '4.34' },
    @{ Cell = 'E33'; Value = '  +3.45%  ' },
    @{ Cell = 'E34'; Value = '  +4.54%  ' },
    @{ Cell = 'D35'; Value = '1.63' },
    @{ Cell = 'E35'; Value = '  +10.04%  ' },
    @{ Cell = 'E36'; Value = '  +3.82%  ' },
    @{ Cell = 'E37'; Value = '  -0.45%  ' },
    @{ Cell = 'D38'; Value = '0.849' },
    @{ Cell = 'E38'; Value = '  -7.17%  ' },
    @{ Cell = 'E39'; Value = '  +0.57%  ' },
    @{ Cell = 'D40'; Value = '17.64' },
    @{ Cell = 'E40'; Value = '  +7.75%  ' },
    @{ Cell = 'E41'; Value = '  +3.08%  ' },
    @{ Cell = 'E42'; Value = '  +3.73%  ' },
    @{ Cell = 'D43'; Value = '0.0666' },
    @{ Cell = 'E43'; Value = '  +2.14%  ' },
    @{ Cell = 'E45'; Value = '  +1.17%  ' },
    @{ Cell = 'D46'; Value = '1.308.56' },
    @{ Cell = 'E46'; Value = '  -3.44%  ' },
    @{ Cell = 'E47'; Value = '  +0.00%  ' },
    @{ Cell = 'D48'; Value = '2.74' },
    @{ Cell = 'E48'; Value = '  -1.53%  ' },
    @{ Cell = 'D49'; Value = '6.61' },
    @{ Cell = 'E49'; Value = '  +1.87%  ' },
    @{ Cell = 'D50'; Value = '12.13' },
    @{ Cell = 'E50'; Value = '  -1.49%  ' },
    @{ Cell = 'D51'; Value = '0.0763' },
    @{ Cell = 'E51'; Value = '  +8.03%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        # Plain-looking number (e.g. "252.62", "5.00") - force Text format
        # first so Excel stores it verbatim instead of converting it to a
        # numeric value (which would also strip trailing zeros).
        $range.NumberFormat = "@"
    }
    $range.Value = $u.Value
}
